$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.269.73"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.493.90"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.11"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.07"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").Value = "3.493.60"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.17"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "4.087.61"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.120"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "3.493.44"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "64.303.40"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.34"
$ws.Range("E18").Value = "  -8.87%  "
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.73"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("E21").Value = "  -5.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.42"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.566"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.633.61"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.32"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.51"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.24"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "3.515.59"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.47"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "162.14"
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.63"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "2.469.76"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("E51").Value = "  -2.08%  "
